# Update the "StatOutput" summary row (row 2) with the new count-query
# results for the Black and Tan Coonhound breed filter:
#   number_of_files  (A2) stays "1"
#   number_of_sample (B2) "2" -> "3"
#   number_of_cases  (C2) "1" -> "2"
#   number_of_study  (D2) "1" -> "2"
#
# These are digit-only strings that must remain stored as TEXT (shared
# string) rather than being auto-coerced to numbers by Excel's normal
# "smart" Value-assignment behaviour. Forcing the NumberFormat to "@"
# (Text) before the assignment keeps them textual; ClearFormats()
# afterwards restores the cells to their original (unstyled) appearance
# without disturbing the underlying text value.
$wb = $excel.ActiveWorkbook
$statOutput = $wb.Worksheets.Item("StatOutput")

$statOutput.Range("B2:D2").NumberFormat = "@"
$statOutput.Range("B2").Value = "3"
$statOutput.Range("C2").Value = "2"
$statOutput.Range("D2").Value = "2"
$statOutput.Range("B2:D2").ClearFormats()

# Update the matching Cypher query text shown on the "StatOutput_Message"
# sheet (row 18) so it reflects the Black and Tan Coonhound breed filter
# instead of the old Akita example.
$statOutputMessage = $wb.Worksheets.Item("StatOutput_Message")
$newCypher = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Black and Tan Coonhound']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$statOutputMessage.Range("A18").Value = $newCypher
